$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record is added at the top of the data block (row 7),
# pushing all the existing records (old rows 7-34) down by one row.
$ws.Rows("7:7").Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A7").Value = 1
$ws.Range("B7").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C7").Value = "Arica y Parinacota"
$ws.Range("D7").Value = "2022-04-08"
$ws.Range("E7").Value = 15
$ws.Range("F7").Value = 100112052
$ws.Range("G7").Value = "Albahaca"
$ws.Range("H7").Value = "Sin especificar"
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 950
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 975
$ws.Range("N7").Value = "$/paquete"
$ws.Range("O7").Value = "Región de Arica y Parinacota"
$ws.Range("P7").Value = 975
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = "Hortaliza"
